$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (row 1)
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# Title-case connector words (de/del/la/las/los/el) in state & municipality names
$ws.Range("B26").Value2 = "Mazapa De Madero"
$ws.Range("B32").Value2 = "San Cristóbal De Las Casas"
$ws.Range("B41").Value2 = "Hidalgo Del Parral"
$ws.Range("A49").Value2 = "Ciudad De México"
$ws.Range("B67").Value2 = "San Juan Del Río"
$ws.Range("A69").Value2 = "Estado De México"
$ws.Range("B69").Value2 = "Almoloya De Alquisiras"
$ws.Range("B70").Value2 = "Almoloya De Juárez"
$ws.Range("B77").Value2 = "Chapa De Mota"
$ws.Range("B81").Value2 = "Ecatepec De Morelos"
$ws.Range("B84").Value2 = "Ixtapan De La Sal"
$ws.Range("B94").Value2 = "Naucalpan De Juárez"
$ws.Range("B98").Value2 = "San Felipe Del Progreso"
$ws.Range("B108").Value2 = "Tenango Del Valle"
$ws.Range("B113").Value2 = "Tlalnepantla De Baz"
$ws.Range("B116").Value2 = "Valle De Chalco Solidaridad"
$ws.Range("B117").Value2 = "Villa De Allende"
$ws.Range("B118").Value2 = "Villa Del Carbón"
$ws.Range("B126").Value2 = "Apaseo El Grande"
$ws.Range("B129").Value2 = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B132").Value2 = "Jaral Del Progreso"
$ws.Range("B141").Value2 = "Valle De Santiago"
$ws.Range("B145").Value2 = "Acapulco De Juárez"
$ws.Range("B146").Value2 = "Ajuchitlán Del Progreso"
$ws.Range("B147").Value2 = "Alcozauca De Guerrero"
$ws.Range("B149").Value2 = "Atoyac De Álvarez"
$ws.Range("B150").Value2 = "Chilpancingo De Los Bravo"
$ws.Range("B151").Value2 = "Coyuca De Benítez"
$ws.Range("B156").Value2 = "Huitzuco De Los Figueroa"
$ws.Range("B157").Value2 = "Iguala De La Independencia"
$ws.Range("B159").Value2 = "Zihuatanejo De Azueta"
$ws.Range("B161").Value2 = "La Unión De Isidoro Montes De Oca"
$ws.Range("B170").Value2 = "Taxco De Alarcón"
$ws.Range("B171").Value2 = "Técpan De Galeana"
$ws.Range("B173").Value2 = "Tepecoacuilco De Trujano"
$ws.Range("B183").Value2 = "Huasca De Ocampo"
$ws.Range("B186").Value2 = "Molango De Escamilla"
$ws.Range("B187").Value2 = "Nopala De Villagrán"
$ws.Range("B188").Value2 = "Pachuca De Soto"
$ws.Range("B194").Value2 = "Autlán De Navarro"
$ws.Range("B197").Value2 = "Encarnación De Díaz"
$ws.Range("B199").Value2 = "Huejuquilla El Alto"
$ws.Range("B201").Value2 = "Lagos De Moreno"
$ws.Range("B204").Value2 = "San Martín De Bolaños"
$ws.Range("B205").Value2 = "Tepatitlán De Morelos"
$ws.Range("B208").Value2 = "Zacoalco De Torres"
$ws.Range("B210").Value2 = "Zapotlán El Grande"
$ws.Range("B220").Value2 = "Coalcomán De Vázquez Pallares"
$ws.Range("B261").Value2 = "Puente De Ixtla"
$ws.Range("B265").Value2 = "Tlaltizapán De Zapata"
$ws.Range("B273").Value2 = "Santa María Del Oro"
$ws.Range("B276").Value2 = "Acatlán De Pérez Figueroa"
$ws.Range("B278").Value2 = "Fresnillo De Trujano"
$ws.Range("B279").Value2 = "Guevea De Humboldt"
$ws.Range("B280").Value2 = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B284").Value2 = "Oaxaca De Juárez"
$ws.Range("B285").Value2 = "Putla Villa De Guerrero"
$ws.Range("B293").Value2 = "San Miguel Del Puerto"
$ws.Range("B299").Value2 = "Santo Domingo De Morelos"
$ws.Range("B301").Value2 = "Villa Sola De Vega"
$ws.Range("B322").Value2 = "Huehuetlán El Chico"
$ws.Range("B324").Value2 = "Ixcamilpa De Guerrero"
$ws.Range("B325").Value2 = "Izúcar De Matamoros"
$ws.Range("B334").Value2 = "San Salvador El Seco"
$ws.Range("B335").Value2 = "San Salvador El Verde"
$ws.Range("B336").Value2 = "Tecali De Herrera"
$ws.Range("B340").Value2 = "Tlacotepec De Benito Juárez"
$ws.Range("B350").Value2 = "Cadereyta De Montes"
$ws.Range("B352").Value2 = "Jalpan De Serra"
$ws.Range("B353").Value2 = "Pinal De Amoles"
$ws.Range("B359").Value2 = "Axtla De Terrazas"
$ws.Range("B361").Value2 = "Ciudad Del Maíz"
$ws.Range("B371").Value2 = "Villa De Reyes"
$ws.Range("B399").Value2 = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B405").Value2 = "Cosamaloapan De Carpio"
$ws.Range("B409").Value2 = "Ignacio De La Llave"
$ws.Range("B410").Value2 = "Ixhuatlán De Madero"
$ws.Range("B414").Value2 = "Juchique De Ferrer"
$ws.Range("B418").Value2 = "Martínez De La Torre"
$ws.Range("B420").Value2 = "Medellín De Bravo"
$ws.Range("B429").Value2 = "Soledad De Doblado"

# Remove trailing metadata/footer rows (444-448) and shrink used range to A1:D442
$ws.Range("A444:A448").EntireRow.Delete()

